$wb = $excel.ActiveWorkbook

# --- FT_indicators sheet: insert new row for IndicatorID 71 (England) ---
$wsFT = $wb.Worksheets.Item("FT_indicators")

# Insert a new row above the old row 19 (111/22001/LA), shifting it down to row 20
$wsFT.Rows("19:19").Insert()

# Copy formatting from an existing data row (row 2) onto the new row 19 so the
# new cells pick up the same style (s="1") used by the rest of the table
$wsFT.Range("A2:B2").Copy()
$wsFT.Range("A19:B19").PasteSpecial(-4122)

$wsFT.Range("A19").Value = 71
$wsFT.Range("B19").Value = 91041
$wsFT.Range("C19").Value = "England"

# Grow Table1 so it covers the newly inserted row
$loFT = $wsFT.ListObjects.Item(1)
$loFT.Resize($wsFT.Range("A1:C20"))

# --- meta_only sheet: remove the now-duplicated 71/91041 row ---
$wsMeta = $wb.Worksheets.Item("meta_only")
$wsMeta.Rows("4:4").Delete()

# --- Window / selection state ---
$wsFT.Activate()
$wsFT.Range("B19").Select() | Out-Null

$wsMeta.Activate()
$wsMeta.Range("H21").Select() | Out-Null
